# Insert a new weekly price record for "Pepino ensalada" (Feria Lagunitas
# de Puerto Montt) above the existing row 54, shifting the rest of the
# table down by one row (old row 54 -> 55, ..., old row 167 -> 168).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 54:167 down to 55:168, leaving a blank row 54 to fill in.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Cells.Item(54, 1).Value = 4
$ws.Cells.Item(54, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(54, 3).Value = "Los Lagos"
$ws.Cells.Item(54, 4).Value = 44498
$ws.Cells.Item(54, 5).Value = 10
$ws.Cells.Item(54, 6).Value = 100112043
$ws.Cells.Item(54, 7).Value = "Pepino ensalada"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 400
$ws.Cells.Item(54, 11).Value = 12000
$ws.Cells.Item(54, 12).Value = 12000
$ws.Cells.Item(54, 13).Value = 12000
$ws.Cells.Item(54, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(54, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value = 200
$ws.Cells.Item(54, 17).Value = 60
$ws.Cells.Item(54, 18).Value = "Hortaliza"
